# Updated cryptos list values (Price / Volume(1h)) per the target commit.
# Cells whose new text looks like a plain number are written with a leading
# apostrophe so Excel stores them as literal text (matching the original
# inline-string cell content) instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '30.370.28' }
    @{ Cell = "E2"; Value = '  -3.03%  ' }
    @{ Cell = "D3"; Value = '1.940.25' }
    @{ Cell = "E3"; Value = '  -2.95%  ' }
    @{ Cell = "D4"; Value = '''1.002' }
    @{ Cell = "E4"; Value = '  -0.42%  ' }
    @{ Cell = "D5"; Value = '''251.73' }
    @{ Cell = "E5"; Value = '  -1.56%  ' }
    @{ Cell = "D6"; Value = '''0.7119' }
    @{ Cell = "E6"; Value = '  -5.34%  ' }
    @{ Cell = "D7"; Value = '''1.003' }
    @{ Cell = "E7"; Value = '  -0.12%  ' }
    @{ Cell = "D8"; Value = '''0.3308' }
    @{ Cell = "E8"; Value = '  -3.54%  ' }
    @{ Cell = "D9"; Value = '''27.30' }
    @{ Cell = "E9"; Value = '  -0.94%  ' }
    @{ Cell = "D10"; Value = '''0.07345' }
    @{ Cell = "E10"; Value = '  +2.59%  ' }
    @{ Cell = "D11"; Value = '''0.8066' }
    @{ Cell = "E11"; Value = '  -3.39%  ' }
    @{ Cell = "D12"; Value = '''0.08081' }
    @{ Cell = "E12"; Value = '  -1.40%  ' }
    @{ Cell = "D13"; Value = '1.936.07' }
    @{ Cell = "E13"; Value = '  -3.48%  ' }
    @{ Cell = "D14"; Value = '''5.501' }
    @{ Cell = "E14"; Value = '  -2.05%  ' }
    @{ Cell = "D15"; Value = '''94.52' }
    @{ Cell = "E15"; Value = '  -6.08%  ' }
    @{ Cell = "E16"; Value = '  -2.86%  ' }
    @{ Cell = "D17"; Value = '30.363.13' }
    @{ Cell = "E17"; Value = '  -3.05%  ' }
    @{ Cell = "D18"; Value = '''0.000008228' }
    @{ Cell = "E18"; Value = '  -1.02%  ' }
    @{ Cell = "D19"; Value = '''253.02' }
    @{ Cell = "E19"; Value = '  -5.86%  ' }
    @{ Cell = "D20"; Value = '''5.836' }
    @{ Cell = "E20"; Value = '  -4.49%  ' }
    @{ Cell = "D21"; Value = '2.190.92' }
    @{ Cell = "E21"; Value = '  -3.33%  ' }
    @{ Cell = "D22"; Value = '''1.003' }
    @{ Cell = "E22"; Value = '  -0.19%  ' }
    @{ Cell = "D23"; Value = '''1.002' }
    @{ Cell = "E23"; Value = '  -0.58%  ' }
    @{ Cell = "D24"; Value = '''7.015' }
    @{ Cell = "E24"; Value = '  -1.27%  ' }
    @{ Cell = "D25"; Value = '''9.730' }
    @{ Cell = "E25"; Value = '  -3.27%  ' }
    @{ Cell = "D26"; Value = '''163.54' }
    @{ Cell = "E26"; Value = '  -0.33%  ' }
    @{ Cell = "D27"; Value = '''2.357' }
    @{ Cell = "E27"; Value = '  -1.19%  ' }
    @{ Cell = "D28"; Value = '''19.34' }
    @{ Cell = "E28"; Value = '  -3.06%  ' }
    @{ Cell = "D29"; Value = '''0.1304' }
    @{ Cell = "E29"; Value = '  -3.10%  ' }
    @{ Cell = "D30"; Value = '''1.575' }
    @{ Cell = "E30"; Value = '  -1.81%  ' }
    @{ Cell = "D31"; Value = '''1.351' }
    @{ Cell = "E31"; Value = '  -2.38%  ' }
    @{ Cell = "D32"; Value = '''4.423' }
    @{ Cell = "E32"; Value = '  -5.35%  ' }
    @{ Cell = "E33"; Value = '  -6.90%  ' }
    @{ Cell = "D34"; Value = '''0.05190' }
    @{ Cell = "E34"; Value = '  -3.49%  ' }
    @{ Cell = "E35"; Value = '  -1.98%  ' }
    @{ Cell = "D36"; Value = '''0.7471' }
    @{ Cell = "E36"; Value = '  -5.47%  ' }
    @{ Cell = "D37"; Value = '''2.753' }
    @{ Cell = "E37"; Value = '  -1.44%  ' }
    @{ Cell = "D38"; Value = '''0.01973' }
    @{ Cell = "D39"; Value = '''2.813' }
    @{ Cell = "E39"; Value = '  -3.42%  ' }
    @{ Cell = "D40"; Value = '''78.96' }
    @{ Cell = "E40"; Value = '  -7.76%  ' }
    @{ Cell = "D41"; Value = '''6.417' }
    @{ Cell = "E41"; Value = '  -6.12%  ' }
    @{ Cell = "D42"; Value = '''0.4531' }
    @{ Cell = "E42"; Value = '  -2.80%  ' }
    @{ Cell = "E43"; Value = '  -5.53%  ' }
    @{ Cell = "D44"; Value = '''0.8488' }
    @{ Cell = "E44"; Value = '  -1.09%  ' }
    @{ Cell = "D45"; Value = '''1.003' }
    @{ Cell = "E45"; Value = '  -0.09%  ' }
    @{ Cell = "D46"; Value = '''101.67' }
    @{ Cell = "E46"; Value = '  -3.46%  ' }
    @{ Cell = "D47"; Value = '''7.467' }
    @{ Cell = "E47"; Value = '  -4.00%  ' }
    @{ Cell = "D48"; Value = '''9.693' }
    @{ Cell = "E48"; Value = '  -4.58%  ' }
    @{ Cell = "D49"; Value = '''36.68' }
    @{ Cell = "E49"; Value = '  -2.43%  ' }
    @{ Cell = "D50"; Value = '''0.4186' }
    @{ Cell = "E50"; Value = '  -3.67%  ' }
    @{ Cell = "D51"; Value = '''0.06040' }
    @{ Cell = "E51"; Value = '  -0.47%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
